$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "891×6=5346" "881×8=7048"
Replace-Text "309×8=2472" "989×5=4945"
Replace-Text "639×4=2556" "273×2=546"
Replace-Text "166×2=332" "894×3=2682"
Replace-Text "698×5=3490" "721×5=3605"
Replace-Text "960×6=5760" "754×2=1508"
Replace-Text "359×2=718" "285×9=2565"
Replace-Text "215×7=1505" "720×6=4320"
Replace-Text "923×5=4615" "359×7=2513"
Replace-Text "937×8=7496" "790×5=3950"
Replace-Text "714×7=4998" "988×2=1976"
Replace-Text "368×8=2944" "131×4=524"
Replace-Text "786×3=2358" "525×4=2100"
Replace-Text "812×9=7308" "897×6=5382"
Replace-Text "257×4=1028" "219×6=1314"
Replace-Text "419×9=3771" "532×3=1596"
Replace-Text "149×4=596" "194×4=776"
Replace-Text "855×7=5985" "216×3=648"
Replace-Text "631×4=2524" "757×3=2271"
Replace-Text "812×6=4872" "925×7=6475"
Replace-Text "305×3=915" "869×8=6952"
Replace-Text "929×8=7432" "930×6=5580"
Replace-Text "189×6=1134" "388×7=2716"
Replace-Text "723×7=5061" "274×5=1370"
Replace-Text "788×8=6304" "220×4=880"
